# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Updates the worker/period/value table (rows 16-29) on the active sheet
# to reflect the refreshed "estado de cuenta" data: a new worker
# (LINDA MARISOL MONTOYA GARCIA) is inserted at the top of the table,
# the remaining workers' rows shift down, and several mora amounts are
# updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Row, DocNumber(C), Name(D), Period(E), ValorMora(F), SalarioBasico(G)
$rows = @(
    @(16, "1047407520", "LINDA MARISOL MONTOYA GARCIA", "2403", 52000,  1432260),
    @(17, "1128055034", "GINA PAOLA SERRANO PRADA",     "2403", 52000,  126000),
    @(18, "1128055034", "GINA PAOLA SERRANO PRADA",     "2507", 138041, 3451020),
    @(19, "1128055034", "GINA PAOLA SERRANO PRADA",     "2506", 138041, 3451020),
    @(20, "1128055034", "GINA PAOLA SERRANO PRADA",     "2505", 138041, 5235900),
    @(21, "1128055034", "GINA PAOLA SERRANO PRADA",     "2504", 138041, 5235900),
    @(22, "1128055034", "GINA PAOLA SERRANO PRADA",     "2503", 138041, 5235900),
    @(23, "1128055034", "GINA PAOLA SERRANO PRADA",     "2502", 138041, 5235900),
    @(24, "9144427",    "JOSE GREGORIO CASTRO MARTINEZ","2507", 46400,  1160000),
    @(25, "9144427",    "JOSE GREGORIO CASTRO MARTINEZ","2506", 46400,  1160000),
    @(26, "9144427",    "JOSE GREGORIO CASTRO MARTINEZ","2505", 46400,  1160000),
    @(27, "9144427",    "JOSE GREGORIO CASTRO MARTINEZ","2504", 46400,  1160000),
    @(28, "9144427",    "JOSE GREGORIO CASTRO MARTINEZ","2503", 46400,  1160000),
    @(29, "9144427",    "JOSE GREGORIO CASTRO MARTINEZ","2502", 46400,  1160000)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = "CC"
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
}
